# "Adding Course Catalogs Page" -- expand the login/account-creation test
# data from 3 rows to 10 rows and refresh the active-sheet/selection state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("loginTest")
$ws2 = $wb.Worksheets.Item("accountCreatonTest")
$ws3 = $wb.Worksheets.Item("profileUpdateTest")

# New login / student-id values used to grow each data sheet from 3 to 10 rows.
$stds   = @("STD01","STD02","STD03","STD04","STD05","STD06","STD07","STD08","STD09","STD10")
$logins = @("CHORME_020516aa","CHORME_020516ab","CHORME_020516ac","CHORME_020516ad","CHORME_020516ae","CHORME_020516af","CHORME_020516ag","CHORME_020516ah","CHORME_020516ai","CHORME_020516aj")

# ---------------------------------------------------------------------------
# accountCreatonTest (sheet2): columns A (student id) and C/F (login) drive
# the shared-string table order, so write them before anything else --
# column A first (STD01..STD10), then C/F (CHORME_020516aa..aj) -- matching
# the order new strings must appear in the rebuilt sharedStrings table.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 10; $i++) {
    $r = 2 + $i
    $ws2.Cells.Item($r, 1).Value = $stds[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $r = 2 + $i
    $ws2.Cells.Item($r, 3).Value = $logins[$i]
    $ws2.Cells.Item($r, 6).Value = $logins[$i]
}

# Remaining text columns for accountCreatonTest rows 2-11 (values repeat the
# same constants already used in rows 2-4, so they reuse existing shared
# strings rather than minting new ones).
for ($i = 0; $i -lt 10; $i++) {
    $r = 2 + $i
    $ws2.Cells.Item($r, 2).Value  = "M"
    $ws2.Cells.Item($r, 4).Value  = "01/01/1980"
    $ws2.Cells.Item($r, 5).Value  = "jcarter.dsi@gmail.com"
    $ws2.Cells.Item($r, 9).Value  = "q"
    $ws2.Cells.Item($r, 10).Value = "a"
}

# Columns G/H are numeric (12345678), stored as plain numbers even though
# the column's display style is a text format. Force numeric storage via a
# General-format round trip, then copy row 4's cell *format only* back on
# top so the style id matches the rest of the column.
for ($i = 2; $i -lt 10; $i++) {
    $r = 2 + $i
    $ws2.Cells.Item($r, 7).Style = "Normal"
    $ws2.Cells.Item($r, 7).Value = 12345678
    $ws2.Cells.Item(4, 7).Copy()
    $ws2.Cells.Item($r, 7).PasteSpecial(-4122)

    $ws2.Cells.Item($r, 8).Style = "Normal"
    $ws2.Cells.Item($r, 8).Value = 12345678
    $ws2.Cells.Item(4, 8).Copy()
    $ws2.Cells.Item($r, 8).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# loginTest (sheet1): column A gets the same CHORME_020516aa..aj logins
# (reusing the shared strings just created above), column B stays the
# constant text "12345678".
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 10; $i++) {
    $r = 2 + $i
    $ws1.Cells.Item($r, 1).Value = $logins[$i]
    $ws1.Cells.Item($r, 2).Value = "12345678"
}

# ---------------------------------------------------------------------------
# Cosmetic / navigation state to match the saved workbook: widen
# accountCreatonTest's column C, move the active tab/selection from
# profileUpdateTest back to loginTest, and refresh each sheet's selection.
# ---------------------------------------------------------------------------
$ws2.Columns.Item(3).ColumnWidth = 17.5

$ws2.Range("B7").Select()
$ws1.Range("B15").Select()
$ws1.Activate()
